$d = $word.ActiveDocument

# 1. Update version number in revision history table: "2.0" -> "1.1"
$d.Content.Find.Execute("2.0", $false, $false, $false, $false, $false, $true, 1, $false, "1.1", 2)
